$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 1-3: summary values become "0M"
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# Row 4
$t.Cell(4,1).Range.Text = "218"

# Row 5
$t.Cell(5,1).Range.Text = "0.00001"

# Row 6
$t.Cell(6,1).Range.Text = "0.00059"

# Row 7
$t.Cell(7,1).Range.Text = "0.00013"

# Row 9
$t.Cell(9,1).Range.Text = "0.00019"

# Row 10
$t.Cell(10,1).Range.Text = "0.00022"

# Row 11
$t.Cell(11,1).Range.Text = "0.00026"

# Row 12
$t.Cell(12,1).Range.Text = "0.03306"

# Rows 44-46: collapse tab-separated values into the single values
# moved from rows 1-3 (now replaced with "0M" above)
$t.Cell(44,1).Range.Text = "99.97"
$t.Cell(45,1).Range.Text = "0.03"
$t.Cell(46,1).Range.Text = "111"
